# Commit: "Created seed project, created Halls, Movies and Screenings seed methods"
# This particular workbook (Events.xlsx) was touched incidentally: the
# "data" sheet gets its proper name "Events", the author's Excel window
# was a different size/position when they saved, and the cell selection
# left on the sheet moved off the data range.

$wb = $excel.ActiveWorkbook

# Rename the worksheet from "data" to "Events"
$ws = $wb.Worksheets.Item("data")
$ws.Name = "Events"

# Match the author's window geometry at save time (restored / moved window)
$win = $excel.ActiveWindow
$win.WindowState = -4143   # xlNormal
$win.Left = -19320
$win.Top = -120
$win.Width = 19440
$win.Height = 15000

# Leave the selection where the author left it before saving
$ws.Range("C34").Select()
